# Auto update Excel log
# Appends newly-captured sensor events to the "PIR" sheet (rows 36-48)
# and the "Proximity" sheet (rows 3-4).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PIR sheet: additional Bathroom motion-sensor readings for 2026-02-06
# ---------------------------------------------------------------------------
$pirSheet = $wb.Worksheets("PIR")

$pirData = @(
    @("2026-02-06", "09:39:28", "09:00", "Bathroom", "No Motion",       "Inactive"),
    @("2026-02-06", "09:39:29", "09:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-06", "09:39:35", "09:00", "Bathroom", "No Motion",       "Inactive"),
    @("2026-02-06", "09:39:37", "09:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-06", "09:39:45", "09:00", "Bathroom", "No Motion",       "Inactive"),
    @("2026-02-06", "09:39:47", "09:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-06", "09:39:55", "09:00", "Bathroom", "No Motion",       "Inactive"),
    @("2026-02-06", "09:39:55", "09:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-06", "09:40:03", "09:00", "Bathroom", "No Motion",       "Inactive"),
    @("2026-02-06", "09:40:07", "09:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-06", "09:40:15", "09:00", "Bathroom", "No Motion",       "Inactive"),
    @("2026-02-06", "09:40:20", "09:00", "Bathroom", "No Motion",       "Inactive"),
    @("2026-02-06", "09:40:25", "09:00", "Bathroom", "No Motion",       "Inactive")
)

$pirStartRow = 36
$pirEndRow = $pirStartRow + $pirData.Count - 1
$pirRange = $pirSheet.Range("A" + $pirStartRow + ":F" + $pirEndRow)
# Force plain-text storage so date/time-looking strings are not
# reinterpreted as Excel dates/times.
$pirRange.NumberFormat = "@"

for ($i = 0; $i -lt $pirData.Count; $i++) {
    $r = $pirStartRow + $i
    $values = $pirData[$i]
    for ($c = 1; $c -le 6; $c++) {
        $pirSheet.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}

# Restore the default (unstyled) cell style now that the text is set.
$pirRange.Style = "Normal"

# ---------------------------------------------------------------------------
# Proximity sheet: Bathroom Door enter/exit events for 2026-02-06
# ---------------------------------------------------------------------------
$proximitySheet = $wb.Worksheets("Proximity")

$proximityData = @(
    @("2026-02-06", "09:39:41", "09:00", "Bathroom Door", "ENTER", "User ENTERED Bathroom"),
    @("2026-02-06", "09:39:45", "09:00", "Bathroom Door", "EXIT",  "User EXITED Bathroom")
)

$proxStartRow = 3
$proxEndRow = $proxStartRow + $proximityData.Count - 1
$proxRange = $proximitySheet.Range("A" + $proxStartRow + ":F" + $proxEndRow)
$proxRange.NumberFormat = "@"

for ($i = 0; $i -lt $proximityData.Count; $i++) {
    $r = $proxStartRow + $i
    $values = $proximityData[$i]
    for ($c = 1; $c -le 6; $c++) {
        $proximitySheet.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}

$proxRange.Style = "Normal"
